$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(63,1).Value = "MotIMCDetuning"
$ws.Cells.Item(63,2).Value = "An experiment after MOT stage. Scan IMC detuning."
$ws.Cells.Item(63,3).Value = "TOP"
$ws.Cells.Item(63,4).Value = "Full"
$ws.Cells.Item(63,5).Value = "None"
$ws.Cells.Item(63,6).Value = 4
$ws.Cells.Item(63,7).Value = "IMCDetuning"
$ws.Cells.Item(63,8).Style = "Normal"
$ws.Cells.Item(63,9).Value = "None"
$ws.Cells.Item(63,10).Value = "LF"
$ws.Cells.Item(63,11).Value = "RandomPolarization"
$ws.Cells.Item(63,12).Value = 8
$ws.Cells.Item(63,13).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(63,14).Value = 30
$ws.Cells.Item(63,15).Value = "LinearFit1D"
